$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 31; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    $cell.Value2 = -$v
}
